# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the scraped data refresh captured by the diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7861
$ws1.Range("F5").Value = 7861
$ws1.Range("F9").Value = 8631
$ws1.Range("F22").Value = 3906
$ws1.Range("F27").Value = 178
$ws1.Range("F29").Value = 5452
$ws1.Range("F31").Value = 71
$ws1.Range("F36").Value = 2279
$ws1.Range("F40").Value = 4527
$ws1.Range("F43").Value = 42
$ws1.Range("F44").Value = 3559

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7861
$ws4.Range("F5").Value = 7861
$ws4.Range("F9").Value = 8631
$ws4.Range("F21").Value = 151
$ws4.Range("F25").Value = 3906
$ws4.Range("F30").Value = 5452
$ws4.Range("F31").Value = 71
$ws4.Range("F35").Value = 72
$ws4.Range("F36").Value = 2279
$ws4.Range("F42").Value = 4527
$ws4.Range("F45").Value = 42
$ws4.Range("F46").Value = 3559
